{"js": "// Replace the 25 three-digit-division answer strings that live in the\n// document's table cells. Each entry is the exact \"old text\" that must be\n// found (in document order) and the \"new text\" that replaces it, mirroring\n// the OOXML diff where only the <w:t> run text inside each cell's single\n// run changes (formatting / paragraph properties stay untouched).\nconst replacements = [\n  [\"971\u00f78=121, 3\", \"734\u00f73=244, 2\"],\n  [\"119\u00f75=23, 4\", \"793\u00f72=396, 1\"],\n  [\"824\u00f72=412, 0\", \"918\u00f74=229, 2\"],\n  [\"283\u00f76=47, 1\", \"544\u00f73=181, 1\"],\n  [\"173\u00f73=57, 2\", \"479\u00f76=79, 5\"],\n  [\"205\u00f74=51, 1\", \"479\u00f76=79, 5\"],\n  [\"607\u00f78=75, 7\", \"748\u00f73=249, 1\"],\n  [\"571\u00f79=63, 4\", \"845\u00f76=140, 5\"],\n  [\"176\u00f73=58, 2\", \"573\u00f74=143, 1\"],\n  [\"377\u00f72=188, 1\", \"481\u00f72=240, 1\"],\n  [\"551\u00f75=110, 1\", \"194\u00f73=64, 2\"],\n  [\"332\u00f72=166, 0\", \"106\u00f76=17, 4\"],\n  [\"947\u00f75=189, 2\", \"730\u00f73=243, 1\"],\n  [\"678\u00f72=339, 0\", \"730\u00f73=243, 1\"],\n  [\"392\u00f73=130, 2\", \"237\u00f76=39, 3\"],\n  [\"860\u00f75=172, 0\", \"411\u00f74=102, 3\"],\n  [\"587\u00f79=65, 2\", \"565\u00f77=80, 5\"],\n  [\"224\u00f78=28, 0\", \"797\u00f74=199, 1\"],\n  [\"245\u00f77=35, 0\", \"996\u00f79=110, 6\"],\n  [\"777\u00f74=194, 1\", \"513\u00f72=256, 1\"],\n  [\"538\u00f77=76, 6\", \"393\u00f79=43, 6\"],\n  [\"679\u00f73=226, 1\", \"711\u00f75=142, 1\"],\n  [\"724\u00f79=80, 4\", \"574\u00f76=95, 4\"],\n  [\"193\u00f76=32, 1\", \"789\u00f75=157, 4\"],\n  [\"388\u00f72=194, 0\", \"855\u00f72=427, 1\"],\n];\n\n// Walk every paragraph in the document body (this recurses into table\n// cells, which is where all of these answers live) in document order and\n// consume the replacement list sequentially. Because every \"old\" value is\n// unique and we only ever look at the *next* unmatched replacement, this\n// is safe even though a couple of \"new\" values repeat.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet nextIdx = 0;\nfor (const paragraph of paragraphs.items) {\n  if (nextIdx >= replacements.length) {\n    break;\n  }\n  const [oldText, newText] = replacements[nextIdx];\n  if (paragraph.text === oldText) {\n    paragraph.insertText(newText, Word.InsertLocation.replace);\n    nextIdx++;\n  }\n}\n\nawait context.sync();\n\nif (nextIdx !== replacements.length) {\n  throw new Error(\n    \"Only matched \" + nextIdx + \" of \" + replacements.length + \" replacements\"\n  );\n}\n", "ps1": "# Replace the 25 three-digit-division answer strings that live in the\n# document's single table. Each pair is the exact current cell text and the\n# text that must replace it; only the <w:t> run text changes in the source\n# OOXML (paragraph/run formatting is untouched), so we overwrite each\n# cell's Range.Text in place rather than rebuilding the cell.\n$replacements = @(\n  @(\"971\u00f78=121, 3\", \"734\u00f73=244, 2\"),\n  @(\"119\u00f75=23, 4\", \"793\u00f72=396, 1\"),\n  @(\"824\u00f72=412, 0\", \"918\u00f74=229, 2\"),\n  @(\"283\u00f76=47, 1\", \"544\u00f73=181, 1\"),\n  @(\"173\u00f73=57, 2\", \"479\u00f76=79, 5\"),\n  @(\"205\u00f74=51, 1\", \"479\u00f76=79, 5\"),\n  @(\"607\u00f78=75, 7\", \"748\u00f73=249, 1\"),\n  @(\"571\u00f79=63, 4\", \"845\u00f76=140, 5\"),\n  @(\"176\u00f73=58, 2\", \"573\u00f74=143, 1\"),\n  @(\"377\u00f72=188, 1\", \"481\u00f72=240, 1\"),\n  @(\"551\u00f75=110, 1\", \"194\u00f73=64, 2\"),\n  @(\"332\u00f72=166, 0\", \"106\u00f76=17, 4\"),\n  @(\"947\u00f75=189, 2\", \"730\u00f73=243, 1\"),\n  @(\"678\u00f72=339, 0\", \"730\u00f73=243, 1\"),\n  @(\"392\u00f73=130, 2\", \"237\u00f76=39, 3\"),\n  @(\"860\u00f75=172, 0\", \"411\u00f74=102, 3\"),\n  @(\"587\u00f79=65, 2\", \"565\u00f77=80, 5\"),\n  @(\"224\u00f78=28, 0\", \"797\u00f74=199, 1\"),\n  @(\"245\u00f77=35, 0\", \"996\u00f79=110, 6\"),\n  @(\"777\u00f74=194, 1\", \"513\u00f72=256, 1\"),\n  @(\"538\u00f77=76, 6\", \"393\u00f79=43, 6\"),\n  @(\"679\u00f73=226, 1\", \"711\u00f75=142, 1\"),\n  @(\"724\u00f79=80, 4\", \"574\u00f76=95, 4\"),\n  @(\"193\u00f76=32, 1\", \"789\u00f75=157, 4\"),\n  @(\"388\u00f72=194, 0\", \"855\u00f72=427, 1\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$used = @()\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n  $used += $false\n}\n$matched = 0\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $range = $cell.Range\n    # Cell text carries a trailing cell-mark (CR + BEL); strip it for comparison.\n    $txt = $range.Text.TrimEnd([char]13, [char]7)\n    for ($i = 0; $i -lt $replacements.Count; $i++) {\n      if ((-not $used[$i]) -and ($txt -eq $replacements[$i][0])) {\n        $range.Text = $replacements[$i][1]\n        $used[$i] = $true\n        $matched = $matched + 1\n        break\n      }\n    }\n  }\n}\n\nWrite-Output \"matched $matched of $($replacements.Count)\"\n"}
